$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 14.21875
$ws.Columns.Item(2).ColumnWidth = 39.77734375
$ws.Columns.Item(3).ColumnWidth = 20.77734375
$ws.Columns.Item(4).ColumnWidth = 22.6640625
$ws.Columns.Item(5).ColumnWidth = 16.88671875

# Default font for whole sheet (Times New Roman 11, centered)
$ws.Cells.Font.Name = "Times New Roman"
$ws.Cells.Font.Size = 11
$ws.Cells.HorizontalAlignment = -4108
$ws.Cells.VerticalAlignment = -4108

# Header values
$ws.Range("A1").Value = "Số hiệu"
$ws.Range("B1").Value = "Tên chứng chỉ"
$ws.Range("C1").Value = "Mã nhân viên"
$ws.Range("D1").Value = "Tên nhân viên"
$ws.Range("E1").Value = "Ngày cấp"

# Header row formatting
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Name = "Times New Roman"
$headerRange.Font.Size = 16
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.RowHeight = 21
$headerRange.Interior.ThemeColor = 7
$headerRange.Interior.TintAndShade = 0.59999389629810485
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Row heights for rows 2-5
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8

# Selection
$ws.Range("D13").Select()
